# Auto-generated PowerShell COM-interop script
# Applies the 'Finita la slide con lo schema del C/S in TCP' edit to slide 4

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Remove the old placeholder reminder textbox (id=21, 'Promemoria inserire il testo nelle varie caselle')
foreach ($sh in @($s.Shapes)) {
    if ($sh.Name -eq "CasellaDiTesto 20") {
        $sh.Delete()
    }
}

# The slide's shape-id allocator never reuses ids that were already handed out
# earlier in the session (matches real PowerPoint). The original authoring session
# created/deleted several scratch shapes while laying out the diagram labels, which
# burned ids 5,6,18,20,22 (21 was already used by the box we just removed, 12 by the
# picture). Reproduce that so the surviving shapes land on the same ids as the target.
function New-ScratchShape($slide) {
    $tmp = $slide.Shapes.AddTextbox(1, 0, 0, 10, 10)
    $tmp.Delete()
}


# id=4 -> "Client"
$sp1 = $s.Shapes.AddTextbox(1, 152.66669291338582, 98.00007874015748, 63.009370078740154, 29.081259842519685)
$sp1.Name = "CasellaDiTesto 3"
$sp1.Fill.Visible = $false
$sp1.TextFrame.WordWrap = $false
$sp1.TextFrame.AutoSize = 1
$sp1.TextFrame.TextRange.Text = "Client"
$sp1.TextFrame.TextRange.Font.Name = "Avenir Next LT Pro Light"
$sp1.TextFrame.TextRange.LanguageID = "it-IT"
New-ScratchShape $s
New-ScratchShape $s

# id=7 -> "Server"
$sp2 = $s.Shapes.AddTextbox(1, 601.0, 98.00007874015748, 66.66984251968503, 29.081259842519685)
$sp2.Name = "CasellaDiTesto 6"
$sp2.Fill.Visible = $false
$sp2.TextFrame.WordWrap = $false
$sp2.TextFrame.AutoSize = 1
$sp2.TextFrame.TextRange.Text = "Server"
$sp2.TextFrame.TextRange.Font.Name = "Avenir Next LT Pro Light"
$sp2.TextFrame.TextRange.LanguageID = "it-IT"

# id=8 -> "Directory 2"
$sp3 = $s.Shapes.AddTextbox(1, 366.1070866141732, 102.84692913385827, 84.46188976377952, 24.234409448818898)
$sp3.Name = "CasellaDiTesto 7"
$sp3.Fill.Visible = $false
$sp3.TextFrame.WordWrap = $false
$sp3.TextFrame.AutoSize = 1
$sp3.TextFrame.TextRange.Text = "Directory 2"
$sp3.TextFrame.TextRange.Font.Name = "Avenir Next LT Pro Light"
$sp3.TextFrame.TextRange.Font.Size = 14
$sp3.TextFrame.TextRange.LanguageID = "it-IT"

# id=9 -> "Directory 1"
$sp4 = $s.Shapes.AddTextbox(1, 449.40094488188976, 208.8467716535433, 84.46188976377952, 24.234409448818898)
$sp4.Name = "CasellaDiTesto 8"
$sp4.Fill.Visible = $false
$sp4.TextFrame.WordWrap = $false
$sp4.TextFrame.AutoSize = 1
$sp4.TextFrame.TextRange.Text = "Directory 1"
$sp4.TextFrame.TextRange.Font.Name = "Avenir Next LT Pro Light"
$sp4.TextFrame.TextRange.Font.Size = 14
$sp4.TextFrame.TextRange.LanguageID = "it-IT"

# id=10 -> "Directory 2"
$sp5 = $s.Shapes.AddTextbox(1, 592.8046456692913, 208.84685039370078, 84.46188976377952, 24.234409448818898)
$sp5.Name = "CasellaDiTesto 9"
$sp5.Fill.Visible = $false
$sp5.TextFrame.WordWrap = $false
$sp5.TextFrame.AutoSize = 1
$sp5.TextFrame.TextRange.Text = "Directory 2"
$sp5.TextFrame.TextRange.Font.Name = "Avenir Next LT Pro Light"
$sp5.TextFrame.TextRange.Font.Size = 14
$sp5.TextFrame.TextRange.LanguageID = "it-IT"

# id=11 -> "Directory N"
$sp6 = $s.Shapes.AddTextbox(1, 732.5870078740157, 208.84685039370078, 87.23874015748031, 24.234409448818898)
$sp6.Name = "CasellaDiTesto 10"
$sp6.Fill.Visible = $false
$sp6.TextFrame.WordWrap = $false
$sp6.TextFrame.AutoSize = 1
$sp6.TextFrame.TextRange.Text = "Directory N"
$sp6.TextFrame.TextRange.Font.Name = "Avenir Next LT Pro Light"
$sp6.TextFrame.TextRange.Font.Size = 14
$sp6.TextFrame.TextRange.LanguageID = "it-IT"

# id=13 -> "Directory 1"
$sp7 = $s.Shapes.AddTextbox(1, 450.2697637795276, 323.9237007874016, 84.46188976377952, 24.234409448818898)
$sp7.Name = "CasellaDiTesto 12"
$sp7.Fill.Visible = $false
$sp7.TextFrame.WordWrap = $false
$sp7.TextFrame.AutoSize = 1
$sp7.TextFrame.TextRange.Text = "Directory 1"
$sp7.TextFrame.TextRange.Font.Name = "Avenir Next LT Pro Light"
$sp7.TextFrame.TextRange.Font.Size = 14
$sp7.TextFrame.TextRange.LanguageID = "it-IT"

# id=14 -> "Directory 2"
$sp8 = $s.Shapes.AddTextbox(1, 592.103937007874, 323.9237007874016, 84.46188976377952, 24.234409448818898)
$sp8.Name = "CasellaDiTesto 13"
$sp8.Fill.Visible = $false
$sp8.TextFrame.WordWrap = $false
$sp8.TextFrame.AutoSize = 1
$sp8.TextFrame.TextRange.Text = "Directory 2"
$sp8.TextFrame.TextRange.Font.Name = "Avenir Next LT Pro Light"
$sp8.TextFrame.TextRange.Font.Size = 14
$sp8.TextFrame.TextRange.LanguageID = "it-IT"

# id=15 -> "File"
$sp9 = $s.Shapes.AddTextbox(1, 757.7048818897638, 324.3109448818898, 37.00291338582677, 24.234409448818898)
$sp9.Name = "CasellaDiTesto 14"
$sp9.Fill.Visible = $false
$sp9.TextFrame.WordWrap = $false
$sp9.TextFrame.AutoSize = 1
$sp9.TextFrame.TextRange.Text = "File"
$sp9.TextFrame.TextRange.Font.Name = "Avenir Next LT Pro Light"
$sp9.TextFrame.TextRange.Font.Size = 14
$sp9.TextFrame.TextRange.LanguageID = "it-IT"

# id=16 -> "DIR"
$sp10 = $s.Shapes.AddTextbox(1, 445.8920472440945, 431.96181102362203, 37.26031496062992, 24.234409448818898)
$sp10.Name = "CasellaDiTesto 15"
$sp10.Fill.Visible = $false
$sp10.TextFrame.WordWrap = $false
$sp10.TextFrame.AutoSize = 1
$sp10.TextFrame.TextRange.Text = "DIR"
$sp10.TextFrame.TextRange.Font.Name = "Avenir Next LT Pro Light"
$sp10.TextFrame.TextRange.Font.Size = 14
$sp10.TextFrame.TextRange.LanguageID = "it-IT"

# id=17 -> "DIR"
$sp11 = $s.Shapes.AddTextbox(1, 561.4113385826772, 431.96181102362203, 37.26031496062992, 24.234409448818898)
$sp11.Name = "CasellaDiTesto 16"
$sp11.Fill.Visible = $false
$sp11.TextFrame.WordWrap = $false
$sp11.TextFrame.AutoSize = 1
$sp11.TextFrame.TextRange.Text = "DIR"
$sp11.TextFrame.TextRange.Font.Name = "Avenir Next LT Pro Light"
$sp11.TextFrame.TextRange.Font.Size = 14
$sp11.TextFrame.TextRange.LanguageID = "it-IT"
New-ScratchShape $s

# id=19 -> "DIR"
$sp12 = $s.Shapes.AddTextbox(1, 667.021811023622, 431.9617322834646, 37.26031496062992, 24.234409448818898)
$sp12.Name = "CasellaDiTesto 18"
$sp12.Fill.Visible = $false
$sp12.TextFrame.WordWrap = $false
$sp12.TextFrame.AutoSize = 1
$sp12.TextFrame.TextRange.Text = "DIR"
$sp12.TextFrame.TextRange.Font.Name = "Avenir Next LT Pro Light"
$sp12.TextFrame.TextRange.Font.Size = 14
$sp12.TextFrame.TextRange.LanguageID = "it-IT"
New-ScratchShape $s
New-ScratchShape $s

# id=23 -> "FILE"
$sp13 = $s.Shapes.AddTextbox(1, 496.5514960629921, 431.9616535433071, 41.55188976377953, 24.234409448818898)
$sp13.Name = "CasellaDiTesto 22"
$sp13.Fill.Visible = $false
$sp13.TextFrame.WordWrap = $false
$sp13.TextFrame.AutoSize = 1
$sp13.TextFrame.TextRange.Text = "FILE"
$sp13.TextFrame.TextRange.Font.Name = "Avenir Next LT Pro Light"
$sp13.TextFrame.TextRange.Font.Size = 14
$sp13.TextFrame.TextRange.LanguageID = "it-IT"

# id=24 -> "FILE"
$sp14 = $s.Shapes.AddTextbox(1, 612.0707874015748, 431.96157480314963, 41.55188976377953, 24.234409448818898)
$sp14.Name = "CasellaDiTesto 23"
$sp14.Fill.Visible = $false
$sp14.TextFrame.WordWrap = $false
$sp14.TextFrame.AutoSize = 1
$sp14.TextFrame.TextRange.Text = "FILE"
$sp14.TextFrame.TextRange.Font.Name = "Avenir Next LT Pro Light"
$sp14.TextFrame.TextRange.Font.Size = 14
$sp14.TextFrame.TextRange.LanguageID = "it-IT"

# id=25 -> "NOMI DEI DIRETTORI E DEI FILE"
$sp15 = $s.Shapes.AddTextbox(1, 221.658031496063, 488.38275590551183, 222.511968503937, 24.234409448818898)
$sp15.Name = "CasellaDiTesto 24"
$sp15.Fill.Visible = $false
$sp15.TextFrame.WordWrap = $false
$sp15.TextFrame.AutoSize = 1
$sp15.TextFrame.TextRange.Text = "NOMI DEI DIRETTORI E DEI FILE"
$sp15.TextFrame.TextRange.Font.Name = "Avenir Next LT Pro Light"
$sp15.TextFrame.TextRange.Font.Size = 14
$sp15.TextFrame.TextRange.LanguageID = "it-IT"
